# Applies the "fixed indentation issue and data_type issue in calc_bioreactor,
# elaborated error statements" edit to Bioreaktor_forPy.xlsx.
#
# The change touches the "Hübsch" worksheet: several input values in the
# fermentation parameter tables are updated or cleared, which ripples through
# the dependent formulas on "Hübsch" and "Input_Array" automatically.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Hübsch")

# --- Row 3: time-step inputs (N3/O3 doubled, P3/Q3 cleared) ---
$ws.Range("N3").Value = 20
$ws.Range("O3").Value = 20
$ws.Range("P3").ClearContents()
$ws.Range("Q3").ClearContents()

# --- Row 7-9: substrate 1 measured values ---
$ws.Range("D7").Value = 0.2
$ws.Range("D8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("G9").ClearContents()

# --- Rows 14-18: substrate 2 feed table ---
$ws.Range("Q14").Value = 5
$ws.Range("R15").Value = 4
$ws.Range("R16").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("R17").ClearContents()
$ws.Range("D18").ClearContents()

# --- Rows 23-25: product feed table ---
$ws.Range("Q23").Value = 10
$ws.Range("Q24").ClearContents()
$ws.Range("Q25").ClearContents()

# --- Sheet selections / active tab ---
$ws2 = $wb.Worksheets.Item("Input_Array")
$ws2.Range("A33").Select()

$ws.Activate()
$ws.Range("G8:G9").Select()
